$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 323, shifting existing rows 323:432 down to 324:433
$ws.Rows("323:323").Insert()

# Populate the new row 323 with the new record
$ws.Range("A323").Value = 3
$ws.Range("B323").Value = "Femacal de La Calera"
$ws.Range("C323").Value = "Coquimbo"
$ws.Range("D323").Value = 44809
$ws.Range("D323").NumberFormat = $ws.Range("D324").NumberFormat
$ws.Range("E323").Value = 5
$ws.Range("F323").Value = 100112031
$ws.Range("G323").Value = "Poroto verde"
$ws.Range("H323").Value = "Magnum"
$ws.Range("I323").Value = "Primera"
$ws.Range("J323").Value = 85
$ws.Range("K323").Value = 36000
$ws.Range("L323").Value = 37000
$ws.Range("M323").Value = 36529
$ws.Range("N323").Value = "`$/malla 25 kilos"
$ws.Range("O323").Value = "Región de Arica y Parinacota"
$ws.Range("P323").Value = 1461
$ws.Range("Q323").Value = 25
$ws.Range("R323").Value = "Hortaliza"
